$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the style (border/font/date numfmt) from the last existing data row (A385)
# down through the new date column cells (A386:A464), matching style index "2".
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New daily rows: row, date-serial(A), nuovi pos.(B), somma mobile 7gg.(C), somma mobile 7gg. per 100mila abitanti(D)
$data = @"
386,44460,0,4,39.72194637537239
387,44461,1,3,29.7914597815293
388,44462,2,5,49.65243296921549
389,44463,0,4,39.72194637537239
390,44464,0,4,39.72194637537239
391,44465,1,4,39.72194637537239
392,44466,4,8,79.44389275074478
393,44467,0,8,79.44389275074478
394,44468,0,7,69.51340615690168
395,44469,3,8,79.44389275074478
396,44470,1,9,89.37437934458789
397,44471,0,9,89.37437934458789
398,44472,0,8,79.44389275074478
399,44473,1,5,49.65243296921549
400,44474,0,5,49.65243296921549
401,44475,2,7,69.51340615690168
402,44476,4,8,79.44389275074478
403,44477,0,7,69.51340615690168
404,44478,0,7,69.51340615690168
405,44479,0,7,69.51340615690168
406,44480,2,8,79.44389275074478
407,44481,0,8,79.44389275074478
408,44482,0,6,59.5829195630586
409,44483,0,2,19.8609731876862
410,44484,0,2,19.8609731876862
411,44485,0,2,19.8609731876862
412,44486,0,2,19.8609731876862
413,44487,1,1,9.930486593843098
414,44488,0,1,9.930486593843098
415,44489,0,1,9.930486593843098
416,44490,0,1,9.930486593843098
417,44491,3,4,39.72194637537239
418,44492,0,4,39.72194637537239
419,44493,0,4,39.72194637537239
420,44494,0,3,29.7914597815293
421,44495,0,3,29.7914597815293
422,44496,0,3,29.7914597815293
423,44497,0,3,29.7914597815293
424,44498,1,1,9.930486593843098
425,44499,0,1,9.930486593843098
426,44500,0,1,9.930486593843098
427,44501,0,1,9.930486593843098
428,44502,0,1,9.930486593843098
429,44503,0,1,9.930486593843098
430,44504,0,1,9.930486593843098
431,44505,0,0,0
432,44506,1,1,9.930486593843098
433,44507,0,1,9.930486593843098
434,44508,1,2,19.8609731876862
435,44509,0,2,19.8609731876862
436,44510,0,2,19.8609731876862
437,44511,2,4,39.72194637537239
438,44512,0,4,39.72194637537239
439,44513,0,3,29.7914597815293
440,44514,1,4,39.72194637537239
441,44515,0,3,29.7914597815293
442,44516,4,7,69.51340615690168
443,44517,1,8,79.44389275074478
444,44518,2,8,79.44389275074478
445,44519,0,8,79.44389275074478
446,44520,1,9,89.37437934458789
447,44521,0,8,79.44389275074478
448,44522,0,8,79.44389275074478
449,44523,0,4,39.72194637537239
450,44524,14,17,168.8182720953327
451,44525,0,15,148.9572989076465
452,44526,2,17,168.8182720953327
453,44527,1,17,168.8182720953327
454,44528,0,17,168.8182720953327
455,44529,0,17,168.8182720953327
456,44530,0,17,168.8182720953327
457,44531,0,3,29.7914597815293
458,44532,0,3,29.7914597815293
459,44533,1,2,19.8609731876862
460,44534,2,3,29.7914597815293
461,44535,0,3,29.7914597815293
462,44536,5,8,79.44389275074478
463,44537,0,8,79.44389275074478
464,44538,1,9,89.37437934458789
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = [double]$parts[1]
    $ws.Cells.Item($r, 2).Value = [double]$parts[2]
    $ws.Cells.Item($r, 3).Value = [double]$parts[3]
    $ws.Cells.Item($r, 4).Value = [double]$parts[4]
}
